# auto commit by win-upload.bat
#
# The author duplicated slide 3 (the circuit-board "Sensor Shield" slide)
# to create a new slide, inserted right after it as the new slide 4, and
# then nudged the picture on the original slide 3 up slightly.

$p = $ppt.ActivePresentation

# EMU <-> point helper (PowerPoint COM positions/sizes are in points;
# OOXML stores EMUs -- 914400 EMU per inch, 72 points per inch).
$emuPerPt = 914400 / 72

# 1) Duplicate slide 3 -> new slide lands immediately after it (slide 4).
$sourceSlide = $p.Slides.Item(3)
$sourceSlide.Duplicate() | Out-Null

# 2) Nudge the picture on the original slide 3 up slightly
#    (y: 450000 EMU -> 436618 EMU; x stays at 2160990 EMU).
$slide3Pic = $p.Slides.Item(3).Shapes.Item(1)
$slide3Pic.Top = 436618 / $emuPerPt
